$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.317.37"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "2.249.43"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'307.64"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").Value = "'96.64"
$ws.Range("E6").Value = "  -1.27%  "

$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").Value = "'35.19"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").Value = "'7.29"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "2.591.97"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").Value = "2.247.50"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").Value = "'0.836"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "'13.65"
$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("D18").Value = "44.121.01"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").Value = "0.0₃0971"
$ws.Range("E19").Value = "  +0.79%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.18"
$ws.Range("E20").Value = "  -6.75%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.39"
$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("D22").Value = "'65.77"

$ws.Range("D23").Value = "'237.95"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  -0.54%  "

$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").Value = "'38.84"
$ws.Range("E27").Value = "  +6.66%  "

$ws.Range("D28").Value = "'9.99"
$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("D30").Value = "'5.94"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").Value = "'20.12"
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").Value = "'152.37"
$ws.Range("E32").Value = "  -2.96%  "

$ws.Range("D33").Value = "'0.0799"
$ws.Range("E33").Value = "  -3.97%  "

$ws.Range("D34").Value = "'3.24"
$ws.Range("E34").Value = "  -0.52%  "

$ws.Range("D35").Value = "'2.61"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("E36").Value = "  +3.00%  "

$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("E38").Value = "  -6.48%  "

$ws.Range("D39").Value = "'3.61"
$ws.Range("E39").Value = "  +2.41%  "

$ws.Range("D40").Value = "'14.64"
$ws.Range("E40").Value = "  -5.96%  "

$ws.Range("D41").Value = "'3.86"
$ws.Range("E41").Value = "  -3.76%  "

$ws.Range("E42").Value = "  -1.87%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").Value = "1.755.81"
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("D45").Value = "'83.19"
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'15.42"
$ws.Range("E47").Value = "  +12.92%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'100.54"
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'4.97"
$ws.Range("E49").Value = "  -2.66%  "

$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.58"
$ws.Range("E51").Value = "  -1.76%  "
